# The draft had a grammar slip ("allowing an process") that Word's
# grammar checker had flagged with proofErr markers around "an". The
# final version corrects this to "allowing a process" and, since the
# corrected text is a single contiguous run of words, Word collapses
# the three runs (and drops the now-unnecessary proofErr markers) into
# one run when the replacement is made.
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$found = $find.Execute(
    "allowing an process",  # FindText
    $true,                  # MatchCase
    $false,                 # MatchWholeWord
    $false,                 # MatchWildcards
    $false,                 # MatchSoundsLike
    $false,                 # MatchAllWordForms
    $true,                  # Forward
    1,                      # Wrap (wdFindContinue)
    $false,                 # Format
    "allowing a process",   # ReplaceWith
    2                       # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Expected text 'allowing an process' was not found in the document."
}
